# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-language sheets (zh-cn,
# de-de) now that the two files have been handed back, updates the
# Overview sheet's status text, and widens a few columns that now hold
# longer content.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3b3dee6d12b650451dd0593ca6b48125add95bc/e2e/"
$file1 = "541917e8-e97b-4d43-9828-337e7105f265.md"
$file2 = "bfb6f0a1-2e70-4009-8702-c5262631f7ea.md"

# ---------------------------------------------------------------------------
# Status changes from "In Translation" to the handed-back status message.
# The same text is shown on the Overview sheet (per-language columns) and
# on each language sheet's own "Status" column, so every one of those
# cells needs to be re-written explicitly.
# ---------------------------------------------------------------------------
$handedBack = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $handedBack
$overview.Range("F2").Value = $handedBack
$overview.Range("E3").Value = $handedBack
$overview.Range("F3").Value = $handedBack

# widen the zh-cn / de-de status columns on the Overview sheet
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

$zh.Range("C2").Value = $handedBack
$zh.Range("C3").Value = $handedBack

# row 2 (541917e8 file): Latest Target File + Latest Handback File
$zh.Range("J2").Value = "541917e8-e97b-4d43-9828-337e7105f265.de1783a1c11e789d3491f3fd54ccbabc9ed919b1.zh-cn.xlf"
# row 3 (bfb6f0a1 file): Latest Target File + Latest Handback File
$zh.Range("J3").Value = "bfb6f0a1-2e70-4009-8702-c5262631f7ea.743a615437b9f9568eeb0da414723f13486b66cf.zh-cn.xlf"

# Latest Handback DateTime placeholder gets a real timestamp now
$zh.Range("K2").Value = "2016-08-31 18:28:15"
$zh.Range("K3").Value = "2016-08-31 18:28:15"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40

$de.Range("C2").Value = $handedBack
$de.Range("C3").Value = $handedBack

$de.Range("J2").Value = "541917e8-e97b-4d43-9828-337e7105f265.de1783a1c11e789d3491f3fd54ccbabc9ed919b1.de-de.xlf"
$de.Range("J3").Value = "bfb6f0a1-2e70-4009-8702-c5262631f7ea.743a615437b9f9568eeb0da414723f13486b66cf.de-de.xlf"

$de.Range("K2").Value = "2016-08-31 18:28:22"
$de.Range("K3").Value = "2016-08-31 18:28:22"

# ---------------------------------------------------------------------------
# Add the new "Latest Target File" hyperlinks (column I) on both sheets,
# pointing at the same source .md file as column A. Rebuild the existing
# hyperlink collections (delete + re-add in reading order) so the saved
# relationship ids come out interleaved the way Excel renumbers them
# (A2, I2, A3, I3) instead of simply appended at the end.
# ---------------------------------------------------------------------------
foreach ($ws in @($zh, $de)) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), ($baseUrl + $file1), "", "", $file1)
    $ws.Hyperlinks.Add($ws.Range("I2"), ($baseUrl + $file1), "", "", $file1)
    $ws.Hyperlinks.Add($ws.Range("A3"), ($baseUrl + $file2), "", "", $file2)
    $ws.Hyperlinks.Add($ws.Range("I3"), ($baseUrl + $file2), "", "", $file2)
}
